$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value of 45179 (2023-09-10)
# for every data row (rows 2 through 480). Update it to 45180 (2023-09-11).
$range = $ws.Range("C2:C480")
$range.Value = 45180
